# "Debugged RDA result tables": the single table in this document reports
# RDA variance-partitioning results for three datasets (All, Particle
# Associated, >= 5 Micron Particles), each with four terms (Latitude,
# Depth, log10(Size_Class), Residual). The "% Variance" column (column 5)
# had incorrect values that need to be corrected to the debugged figures.
#
# Row 1 is the header ("Dataset", "Term", ... "% Variance", ...); rows
# 2-13 are the twelve data rows, addressed here by their 1-based table
# row/column index so each cell is targeted unambiguously. (Several of
# the old percentages, e.g. "1.4% ", repeat across different rows, and a
# Range-scoped Find.Execute in this runtime is not confined to the
# supplied range, so per-cell Find/Replace is unsafe here -- direct
# Cell.Range.Text assignment is used instead, which only ever touches the
# addressed cell.)

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$col = 5

$updates = @(
    @{ Row = 2;  New = "7.7% " }   # All / Latitude                          : 3.3%  -> 7.7%
    @{ Row = 3;  New = "3.3% " }   # All / Depth                             : 1.4%  -> 3.3%
    @{ Row = 4;  New = "18.0%" }   # All / log10(Size_Class)                 : 7.6%  -> 18.0%
    @{ Row = 5;  New = "71.1%" }   # All / Residual                          : 30.1% -> 71.1%
    @{ Row = 6;  New = "8.4% " }   # Particle Associated / Latitude          : 3.0%  -> 8.4%
    @{ Row = 7;  New = "3.7% " }   # Particle Associated / Depth             : 1.4%  -> 3.7%
    @{ Row = 8;  New = "27.1%" }   # Particle Associated / log10(Size_Class) : 9.9%  -> 27.1%
    @{ Row = 9;  New = "60.9%" }   # Particle Associated / Residual          : 22.2% -> 60.9%
    @{ Row = 10; New = "9.9% " }   # >= 5 Micron Particles / Latitude        : 2.1%  -> 9.9%
    @{ Row = 11; New = "4.8% " }   # >= 5 Micron Particles / Depth           : 1.0%  -> 4.8%
    @{ Row = 12; New = "22.0%" }   # >= 5 Micron Particles / log10(Size_Class): 4.7% -> 22.0%
    @{ Row = 13; New = "63.3%" }   # >= 5 Micron Particles / Residual        : 13.4% -> 63.3%
)

foreach ($u in $updates) {
    $t.Cell($u.Row, $col).Range.Text = $u.New
}
